$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date: regenerated publish timestamp
$meta.Range("B8").Value = '2025-06-13T15:45:04+00:00'

# FHIR Version: 4.3.0 -> 4.0.1
$meta.Range("B15").Value = '4.0.1'

# ---------------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Row 2 (Attachment) - Constraint(s): drop the "unless an empty Parameters
# resource ... or $this is Parameters" clause from ele-1
$els.Range("AJ2").Value = 'att-1:If the Attachment has data, it SHALL have a contentType {data.empty() or contentType.exists()}' + [char]10 + 'ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}must-be-drs-uri:attachment.url must start with ^drs://. A drs:// hostname-based URI, as defined in the DRS documentation, that tells clients how to access this object. The intent of this field is to make DRS objects self-contained, and therefore easier for clients to store and pass around.  For example, if you arrive at this DRS JSON by resolving a compact identifier-based DRS URI, the `self_uri` presents you with a hostname and properly encoded DRS ID for use in subsequent `access` endpoint calls. {$this.url.matches(''^drs://.*'')}'

# Row 3 (Attachment.id) - Type(s): id -> string
$els.Range("K3").Value = 'string' + [char]10

# Row 5 (Attachment.contentType) - Binding Description / Binding Value Set
$els.Range("Y5").Value = 'The mime type of an attachment. Any valid mime type is allowed.'
$els.Range("Z5").Value = 'http://hl7.org/fhir/ValueSet/mimetypes|4.0.1'

# Row 6 (Attachment.language) - Binding Description
$els.Range("Y6").Value = 'A human language.'

# Row 10 (Attachment.hash) - Comments: R4B -> R4 link
$els.Range("N10").Value = 'The hash is calculated on the data prior to base64 encoding, if the data is based64 encoded. The hash is not intended to support digital signatures. Where protection against malicious threats a digital signature should be considered, see [Provenance.signature](http://hl7.org/fhir/R4/provenance-definitions.html#Provenance.signature) for mechanism to protect a resource with a digital signature.'
